# The "Fixed Assets" sheet is already the active sheet/tab in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new fixed-asset entry (row 2): an active "house" of type "residence"
# acquired in 2020.
$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "house"
$ws.Range("C2").Value = "residence"
$ws.Range("D2").Value = 2020
$ws.Range("E2").Value = 400000
$ws.Range("F2").Value = 600000
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 2090
$ws.Range("I2").Value = 5

# Move the selection to A3, as left by the author after entering the new row.
$ws.Range("A3").Select()
